$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Roscar cigüeñal" -> "Ajustar vástagos a guías", split across
# three runs with identical (bold, sz 21) formatting, matching the commit's
# diff exactly (3 separate <w:r> elements instead of one).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute('Roscar cigüeñal', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found1) {
    $xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p w:rsidR="00A40ABF" w:rsidRDefault="00A40ABF" w:rsidP="00A40ABF">' +
        '<w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:line="220" w:lineRule="exact"/><w:ind w:left="37"/><w:rPr><w:b/><w:sz w:val="21"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="21"/></w:rPr><w:t>Ajustar v</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="21"/></w:rPr><w:t>á</w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:sz w:val="21"/></w:rPr><w:t>stagos a guías</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Change 2: "Rectificar pista reten cigüeñal" -> "Rectificar válvulas"
# (plain single-run text swap).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute('Rectificar pista reten cigüeñal', $true, $false, $false, $false, $false, $true, 1, $false, 'Rectificar válvulas', 2)

# ---------------------------------------------------------------------------
# Change 3: the three runs "${v" + "25" + "}" collapse into a single run
# "${v25}". A plain Find/Replace over the same text merges it back into one
# run automatically.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute('${v25}', $true, $false, $false, $false, $false, $true, 1, $false, '${v25}', 2)
